$d = $word.ActiveDocument

# The document's headers/footers each hold one inline picture (logo).
# This edit simply renames those pictures (the "name" shown in Word's
# Selection Pane / the wp:docPr + pic:cNvPr "name" attribute) -
# swapping image1.jpg <-> image2.jpg for the BTEC logo, and
# image2.png -> image1.png for the two Pearson logos - with no other
# visible content change.
#
# wdHeaderFooterIndex: 1 = primary, 2 = first-page, 3 = even-pages.
$hfIndexes = @(1, 2, 3)

for ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {
    $sec = $d.Sections.Item($secIdx)

    foreach ($hfIdx in $hfIndexes) {
        $header = $sec.Headers.Item($hfIdx)
        if ($header.Exists) {
            $count = $header.Range.InlineShapes.Count
            for ($i = 1; $i -le $count; $i++) {
                $shape = $header.Range.InlineShapes.Item($i)
                $desc = $shape.AlternativeText

                $newName = $null
                if ($desc -eq "BTec_Logo-Orange") {
                    $newName = "image2.jpg"
                }

                if ($newName -ne $null) {
                    # Re-select the shape's own range right before renaming -
                    # some stories need the shape freshly (re)anchored via
                    # Selection before InlineShape.Name can be written.
                    $shape.Range.Select()
                    $word.Selection.InlineShapes.Item(1).Name = $newName
                }
            }
        }

        $footer = $sec.Footers.Item($hfIdx)
        if ($footer.Exists) {
            $count = $footer.Range.InlineShapes.Count
            for ($i = 1; $i -le $count; $i++) {
                $shape = $footer.Range.InlineShapes.Item($i)
                $desc = $shape.AlternativeText

                $newName = $null
                if ($desc -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $newName = "image1.png"
                }

                if ($newName -ne $null) {
                    $shape.Range.Select()
                    $word.Selection.InlineShapes.Item(1).Name = $newName
                }
            }
        }
    }
}

Write-Output "Renamed logo inline shapes."
